$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing Round #2")

# --- Fill in the second tester's (Addison Boyer, columns H:J) scores for
#     Testing Round #2, mirroring the pattern already present in columns M:O.
$scores = @{
    5  = @(5, 5, 5)
    6  = @(5, 5, 5)
    7  = @(5, 5, 5)
    8  = @(5, 5, 5)
    9  = @(5, 3, 5)
    12 = @(5, 5, 5)
    13 = @(5, 5, 5)
    14 = @(5, 5, 5)
    17 = @(5, 5, 5)
    18 = @(5, 3, 5)
    19 = @(5, 5, 5)
    20 = @(5, 5, 5)
    21 = @(5, 5, 5)
}

foreach ($row in $scores.Keys) {
    $vals = $scores[$row]
    $ws.Cells.Item($row, 8).Value = $vals[0]
    $ws.Cells.Item($row, 9).Value = $vals[1]
    $ws.Cells.Item($row, 10).Value = $vals[2]
}

# --- Update Mark Lannen's notes (C24) with the revised wording (adds two
#     trailing spaces vs. the original text).
$ws.Range("C24").Value = "Overall the testing went well, and most particpants were able to complete the tasks relatively easily. The main pain points were users not being able to find the home button and not understanding what was meant when they were asked to 'report' on what they were seeing. Also, the live coding section was too large and users often didn't realize that they needed to scroll down past that section to find the comments. All of these issues have been addressed in the updated prototype and script that will be used for the second round of testing.  "

# --- Add Addison Boyer's round-2 notes (H24), matching the same merged
#     layout/formatting already used for C24:F24.
$ws.Range("C24:F24").Copy()
$ws.Range("H24:K24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H24:K24").Merge()
$ws.Range("H24").Value = "Overall testing seemed to go smoother than last time.  It appeared that users were able to recognize the home button easier.  One user mentioned that the navigation menu should be condensed into a single component, and shouldn't have an over arrow.  One testing participant noticed a typo in the static text in our prototype.  The comments section seemed to be more recognizable when condensed closer to the top of the page.  One user mentioned that it was difficult to tell whether there was content below or not on small screens.  After clicking didn't realize content had popped up.  One user mentioned that instead of having next on the last step of user account creation, have finish."

$wb.Save()
